$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
  "2024-08-28",
  "2024-08-29",
  "2024-08-30",
  "2024-09-02",
  "2024-09-03",
  "2024-09-04",
  "2024-09-05",
  "2024-09-06",
  "2024-09-09",
  "2024-09-10",
  "2024-09-11",
  "2024-09-12",
  "2024-09-13",
  "2024-09-16",
  "2024-09-17",
  "2024-09-18",
  "2024-09-19",
  "2024-09-20",
  "2024-09-23",
  "2024-09-24",
  "2024-09-25"
)

$cvals = @(
  2200.75,
  2193.75,
  2240.199951171875,
  2232.75,
  2240.25,
  2277.25,
  2290.199951171875,
  2256.5,
  2216.800048828125,
  2222.550048828125,
  2209.39990234375,
  2247.5,
  2256.449951171875,
  2251.85009765625,
  2270.39990234375,
  2224.949951171875,
  2171.89990234375,
  2151.699951171875,
  2182.25,
  2215.75,
  2221.10009765625
)
$dvals = @(
  1539.5,
  1499.150024414062,
  1537.550048828125,
  1537.550048828125,
  1530.599975585938,
  1556.550048828125,
  1555.75,
  1559.900024414062,
  1546.25,
  1545.550048828125,
  1591.949951171875,
  1592.849975585938,
  1582.5,
  1577.75,
  1561.699951171875,
  1543.050048828125,
  1515.050048828125,
  1481.099975585938,
  1440.400024414062,
  1414.25,
  1416.400024414062
)
$evals = @(
  1707.449951171875,
  1691.300048828125,
  1731.75,
  1687.900024414062,
  1687.5,
  1686.550048828125,
  1709.449951171875,
  1702.699951171875,
  1704.199951171875,
  1727.849975585938,
  1725.650024414062,
  1747.949951171875,
  1753.699951171875,
  1741.449951171875,
  1713,
  1646.050048828125,
  1649.800048828125,
  1636.75,
  1712.449951171875,
  1697.5,
  1689.199951171875
)
$fvals = @(
  1969.050048828125,
  1961.150024414062,
  1953.800048828125,
  1970.599975585938,
  1924.650024414062,
  1924.650024414062,
  1933.599975585938,
  1928.400024414062,
  1937.099975585938,
  1912.150024414062,
  1867.75,
  1883.349975585938,
  1923.300048828125,
  1900.949951171875,
  1875.599975585938,
  1857,
  1886.5,
  1897.25,
  1952,
  1944.349975585938,
  1909.550048828125
)
$gvals = @(
  1138.300048828125,
  1132.050048828125,
  1127.900024414062,
  1111.550048828125,
  1114,
  1127.900024414062,
  1115.150024414062,
  1100,
  1104.150024414062,
  1113.199951171875,
  1112.599975585938,
  1120.099975585938,
  1118.550048828125,
  1115.849975585938,
  1110.949951171875,
  1079.949951171875,
  1054.449951171875,
  1054.599975585938,
  1055.25,
  1051.550048828125,
  1063.449951171875
)
$hvals = @(
  8555.050048828125,
  8477.400146484375,
  8591.200073242188,
  8540.35009765625,
  8497,
  8572.900146484375,
  8604.14990234375,
  8547.5,
  8508.5,
  8521.300048828125,
  8507.349853515625,
  8591.749877929688,
  8634.5,
  8587.849975585938,
  8531.649780273438,
  8351,
  8277.699951171875,
  8221.39990234375,
  8342.349975585938,
  8323.400024414062,
  8299.700073242188
)
$ivals = @(
  0,
  -0.009076498898377167,
  0.01342391827581785,
  -0.005918844300264037,
  -0.005075915759957742,
  0.008932581674046723,
  0.00364517903223101,
  -0.006584020848860234,
  -0.004562737642585551,
  0.001504383713712758,
  -0.001637097066476197,
  0.009920836202496664,
  0.004975717714982386,
  -0.005402747630327465,
  -0.006544151967287428,
  -0.02117407358786916,
  -0.00877739777608969,
  -0.0068014121265841,
  0.01471161537923817,
  -0.002271536344954651,
  -0.002847388219040138
)
$jvals = @(
  195.3432372519036,
  193.5702045741812,
  196.1686751810184,
  195.0075833360329,
  194.0177412704662,
  195.7508205905787,
  196.4643673773375,
  195.170841886467,
  194.2803285394565,
  194.572600701606,
  194.2540664677808,
  196.1812292428765,
  197.1573716605673,
  196.0921801380266,
  194.8089231116067,
  190.684024638068,
  189.010315104274,
  187.7247780550743,
  190.4865127869734,
  190.0538157499541,
  189.512658754004
)

$startRow = 597
for ($i = 0; $i -lt 21; $i++) {
    $r = $startRow + $i
    $acell = $ws.Range("A$r")
    $acell.NumberFormat = "@"
    $acell.Value = $dates[$i]
    $acell.ClearFormats()
    $ws.Range("C$r").Value = $cvals[$i]
    $ws.Range("D$r").Value = $dvals[$i]
    $ws.Range("E$r").Value = $evals[$i]
    $ws.Range("F$r").Value = $fvals[$i]
    $ws.Range("G$r").Value = $gvals[$i]
    $ws.Range("H$r").Value = $hvals[$i]
    $ws.Range("I$r").Value = $ivals[$i]
    $ws.Range("J$r").Value = $jvals[$i]
}

Write-Output "Inserted rows 597-617"